$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14/15 swap: text columns B (Coin) and C (Link) ---
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# --- Numeric-looking text columns D (Price) and E (Volume(1h)) ---
# Force text NumberFormat first so Excel does not coerce these
# numeric-looking strings (e.g. "1.00", "4.93") into real numbers,
# then restore the default "Normal" style so no stray style index
# attribute is left on the cell (matches the source formatting).
$deCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","E8","D9","E9","D10","E10","E11","D12","E12","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","E27","D28","E28","D29","E29","D30","E30","D31","E31","E32","E33","E34","E35","D36","E36","E37","D39","E39","E40","E41","E42","E43","E44","D45","E45","D46","E46","D47","E47","D48","E48","E49","E50","D51","E51")
foreach ($cellRef in $deCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.230.50"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.490.26"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "584.55"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "172.42"
$ws.Range("E6").Value = "  +3.15%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "2.490.40"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "0.136"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "4.93"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "2.919.63"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "25.47"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "67.071.08"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "2.477.62"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  -5.00%  "
$ws.Range("D21").Value = "350.51"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "4.01"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "68.69"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "2.619.52"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").Value = "0.0₃0903"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").Value = "509.69"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "160.08"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  -6.90%  "
$ws.Range("D39").Value = "18.25"
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "2.36"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("D46").Value = "38.79"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "142.86"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "3.45"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("D51").Value = "0.0730"
$ws.Range("E51").Value = "  -0.72%  "

foreach ($cellRef in $deCells) {
    $ws.Range($cellRef).Style = "Normal"
}
